$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and update the header label in I1) to reflect the new "through" date
$ws.Name = "Through 2022-06-05"
$ws.Range("I1").Value = "2022 (through 06-05)"

# Update June row (row 7) carjacking count for 2022
$ws.Range("I7").Value = 14

# Update Total row (row 14) carjacking count for 2022
$ws.Range("I14").Value = 678
